{"js": "// Manual del Juego: flesh out Gudi's backstory paragraph with the Masamune\n// sword blurb, fix the \"ladea\" typo, and move the \"_GoBack\" bookmark so it\n// sits right after \"Masamune\" (the new edit point) instead of after the\n// boss picture.\n\nconst doc = context.document;\nconst body = doc.body;\n\n// The \"_GoBack\" bookmark currently lives right after the boss picture -\n// remove it from there; it gets re-added at the new edit point below.\ndoc.deleteBookmark(\"_GoBack\");\n\n// Fix the \"ladea\" typo (\"por haber ido de la ladea\" -> \"por haberse ido de\n// la aldea\") and append the new sentence about the Masamune sword. Using\n// search + insertText(\"Replace\") keeps the surrounding Calisto MT / 24pt\n// run formatting intact.\nconst oldText = \"por haber ido de la ladea para entrenar.\";\nconst newText =\n  \"por haberse ido de la aldea para entrenar. Gudi porta la legendaria espada Masamune, que le permite absorber las almas de quienes asesina, y revitalizarse.\";\n\nconst hits = body.search(oldText, { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error(\"Could not find the target sentence to replace.\");\n}\nhits.items[0].insertText(newText, \"Replace\");\nawait context.sync();\n\n// Re-anchor the \"_GoBack\" bookmark right after \"Masamune\" (before the new\n// \", que le permite...\" clause), matching where Word drops it after an edit.\nconst markHits = body.search(\"Gudi porta la legendaria espada Masamune\", {\n  matchCase: true\n});\nmarkHits.load(\"items\");\nawait context.sync();\n\nconst afterMasamune = markHits.items[0].getRange(\"End\");\nafterMasamune.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Manual del Juego: flesh out Gudi's backstory paragraph with the Masamune\n# sword blurb, fix the \"ladea\" typo, and move the \"_GoBack\" bookmark so it\n# sits right after \"Masamune\" (where the new text was inserted) instead of\n# after the boss picture.\n\n$d = $word.ActiveDocument\n\n# The \"_GoBack\" bookmark currently lives right after the boss picture -\n# remove it from there; it gets re-added at the new edit point below.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Fix the \"ladea\" typo (\"por haber ido de la ladea\" -> \"por haberse ido de\n# la aldea\") and append the new sentence about the Masamune sword in one\n# Find/Replace so the whole run keeps its original Calisto MT formatting.\n$oldText = \"por haber ido de la ladea para entrenar.\"\n$newText = \"por haberse ido de la aldea para entrenar. Gudi porta la legendaria espada Masamune, que le permite absorber las almas de quienes asesina, y revitalizarse.\"\n\n$rng = $d.Content\n$found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\nif (-not $found) {\n    throw \"Could not find the target sentence to replace.\"\n}\n\n# Re-anchor the \"_GoBack\" bookmark right after \"Masamune\" (before the new\n# \", que le permite...\" clause), matching where Word drops it after an edit.\n$markRng = $d.Content\n$markRng.Find.Execute(\"Gudi porta la legendaria espada Masamune\") | Out-Null\n$pos = $markRng.End\n$insertionPoint = $d.Range($pos, $pos)\n$d.Bookmarks.Add(\"_GoBack\", $insertionPoint)\n"}
